$d = $word.ActiveDocument

# Locate the word "Proposal" inside the "Project Proposal " heading paragraph
# and collapse the document range to exactly that word (Find.Execute leaves
# the supplied range positioned on the match).
$rng = $d.Content
$rng.Find.Execute("Proposal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Add the reviewer comment anchored to that "Proposal" range.
$comment = $d.Comments.Add($rng, "Hi Spendylove, I think this project proposal sounds fantastic. Great idea.")
$comment.Author = "Chou, Michael"
$comment.Initial = "MC"
